$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1 and name it "Terms"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Terms"

# Attribution / license text block
$ws2.Range("A1").Value = "This dataset on 'SleepData.xlsx' is hypothetical and was generated"
$ws2.Range("A2").Value = "by Paolo G. Hilado (Github: Dcroix) for training purposes on Basic Statistics . Considering"
$ws2.Range("A3").Value = "that most of the values generated by this dataset use randomization, "
$ws2.Range("A4").Value = "in such a rare case that it resembles any existing dataset, it is purely "
$ws2.Range("A5").Value = "coincidental. It is distributed under "
$ws2.Range("A6").Value = " Creative Commons Attribution-NoDerivatives 4.0 International Public License."

# Turn the last line into a hyperlink (applies the built-in Hyperlink style automatically)
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://creativecommons.org/licenses/by-nd/4.0/")

# Widen column A so the long text is readable, and select A11 like the source file
$ws2.Columns.Item(1).ColumnWidth = 72.66
[void]$ws2.Range("A11").Select()

# Make the new "Terms" sheet the active tab
$ws2.Activate()
